$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 11 de Octubre de 2020 a las 07:29"

# Row 27 - Israel
$ws.Range("B27").Value = 290003
$ws.Range("C27").Value = 128
$ws.Range("D27").Value = 225926
$ws.Range("E27").Value = 62136

# Row 59 - Uzbekistan
$ws.Range("B59").Value = 60894
$ws.Range("C59").Value = 118
$ws.Range("E59").Value = 2687
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 503

# Row 82 - Australia
$ws.Range("B82").Value = 27265
$ws.Range("C82").Value = 21
$ws.Range("D82").Value = 24998
$ws.Range("E82").Value = 1369

# Row 186 - Mongolia
$ws.Range("B186").Value = 316
$ws.Range("C186").Value = 1
$ws.Range("E186").Value = 6

# Row 188 - Camboya
$ws.Range("D188").Value = 278
$ws.Range("E188").Value = 5
